$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update J62 from 10.8 to 10.7
$ws.Range("J62").Value = 10.7

# Add new row 63 (quarter 01-04-2021) with updated data.
# A63 must be stored as literal text "01-04-2021" (a shared string), not an
# auto-converted date serial. Writing it as a formula that evaluates to the
# text, then copying/pasting as values, avoids Excel's date auto-recognition
# and keeps the cell on the default (unstyled) format - matching how the
# other date-label cells in column A are stored.
$ws.Range("A63").Formula = "=""01-04-2021"""
$ws.Range("A63").Copy()
$ws.Range("A63").PasteSpecial(-4163)

$ws.Range("B63").Value = 33.1
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 15.2
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 17.5
$ws.Range("G63").Value = 0.3
$ws.Range("H63").Value = 33.6
$ws.Range("I63").Value = 14.7
$ws.Range("J63").Value = 14.5
$ws.Range("K63").Value = 0.2
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("N63").Value = 2.6
$ws.Range("O63").Value = 14.8
$ws.Range("P63").Value = 1.5
$ws.Range("Q63").Value = -0.6
$ws.Range("R63").Value = 0.4
$ws.Range("S63").Value = 13.5
